# Targon race: add three new units (Bettle, Acid Beetle, Fire Beetle) on
# rows 20-22, give "Toxic Spores" (row 19) a second ability ("Charge", added
# before the existing "Blast"), tag the remaining blank Targon rows with the
# race name, and move the window selection to J25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unit names first (matches original authoring order of new shared
#     strings: Bettle, Acid Beetle, Fire Beetle, Assault spikes, Magma jet,
#     Charge). ---
$ws.Range("C20").Value = "Bettle"
$ws.Range("C21").Value = "Acid Beetle"
$ws.Range("C22").Value = "Fire Beetle"

# --- Row 20: Bettle ---
$ws.Range("B20").Value = "Targon"
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 200
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = "g"
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = "h"
$ws.Range("M20").Value = "b"
$ws.Range("N20").Value = 10
$ws.Range("O20").Value = 10
$ws.Range("P20").Value = 10
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = 6
$ws.Range("S20").Value = 3
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("W20").Value = 0
$ws.Range("X20").Value = 0
$ws.Range("Y20").Value = 0
$ws.Range("Z20").Value = 0
$ws.Range("AA20").Value = 0

# --- Row 21: Acid Beetle ---
$ws.Range("B21").Value = "Targon"
$ws.Range("D21").Value = 225
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 200
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 4
$ws.Range("J21").Value = "g"
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = "h"
$ws.Range("M21").Value = "b"
$ws.Range("N21").Value = 10
$ws.Range("O21").Value = 15
$ws.Range("P21").Value = 10
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = 6
$ws.Range("S21").Value = 4
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0
$ws.Range("X21").Value = 0
$ws.Range("Y21").Value = 0
$ws.Range("Z21").Value = 0
$ws.Range("AA21").Value = 0

# --- Row 22: Fire Beetle ---
$ws.Range("B22").Value = "Targon"
$ws.Range("D22").Value = 225
$ws.Range("E22").Value = 40
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 180
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = "g"
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = "h"
$ws.Range("M22").Value = "b"
$ws.Range("N22").Value = 12
$ws.Range("O22").Value = 10
$ws.Range("P22").Value = 12
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 6
$ws.Range("S22").Value = 6
$ws.Range("T22").Value = 0
$ws.Range("U22").Value = 0
$ws.Range("V22").Value = 0
$ws.Range("W22").Value = 0
$ws.Range("X22").Value = 0
$ws.Range("Y22").Value = 0
$ws.Range("Z22").Value = 0
$ws.Range("AA22").Value = 0

# --- Ability columns (AB), filled in after the main grid, row 20 then
#     row 22 -- row 21's "Acid Beetle" has no special ability listed. ---
$ws.Range("AB20").Value = "Assault spikes"
$ws.Range("AB22").Value = "Magma jet"

# --- Row 19 (Toxic Spores): insert "Charge" as the first ability, and push
#     the existing "Blast" ability into the next column (AC19). ---
$existingAbility = $ws.Range("AB19").Value2
$ws.Range("AC19").Value = $existingAbility
$ws.Range("AB19").Value = "Charge"

# --- Rows 23-26: remaining blank Targon rows just get the race label ---
$ws.Range("B23").Value = "Targon"
$ws.Range("B24").Value = "Targon"
$ws.Range("B25").Value = "Targon"
$ws.Range("B26").Value = "Targon"

# --- Move the active selection / view ---
$ws.Range("J25").Select()
